$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Grant Williams", "PF,C", "Charlotte Hornets"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
